$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (new shared strings) to row 1
$ws.Range("D1").Value = "ORG_POR_IDENOLD"
$ws.Range("E1").Value = "ORG_POR_IDENNEW"
$ws.Range("F1").Value = "ORG_POR_STATUS"

# Move the selection to match the edited workbook (user ended on F4)
[void]$ws.Range("F4").Select()
